$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.125.50'
$ws.Range('E2').Value = '  -0.57%  '

# Row 3
$ws.Range('D3').Value = '3.465.33'
$ws.Range('E3').Value = '  -3.45%  '

# Row 4
$ws.Range('E4').Value = '  +0.50%  '

# Row 5
$ws.Range('D5').Value = '''193.63'
$ws.Range('E5').Value = '  -1.81%  '

# Row 6
$ws.Range('D6').Value = '''545.32'
$ws.Range('E6').Value = '  -2.07%  '

# Row 7
$ws.Range('D7').Value = '''0.626'
$ws.Range('E7').Value = '  +2.90%  '

# Row 8
$ws.Range('D8').Value = '3.467.52'
$ws.Range('E8').Value = '  -3.20%  '

# Row 9
$ws.Range('E9').Value = '  -0.03%  '

# Row 10
$ws.Range('D10').Value = '''0.646'
$ws.Range('E10').Value = '  -2.80%  '

# Row 11
$ws.Range('D11').Value = '''59.62'
$ws.Range('E11').Value = '  +6.30%  '

# Row 12
$ws.Range('D12').Value = '''0.140'
$ws.Range('E12').Value = '  -6.93%  '

# Row 13
$ws.Range('D13').Value = '''0.0000261'
$ws.Range('E13').Value = '  -8.73%  '

# Row 14
$ws.Range('D14').Value = '''9.61'
$ws.Range('E14').Value = '  -2.70%  '

# Row 15
$ws.Range('D15').Value = '4.054.19'
$ws.Range('E15').Value = '  -2.51%  '

# Row 16
$ws.Range('D16').Value = '3.493.43'
$ws.Range('E16').Value = '  -2.93%  '

# Row 17
$ws.Range('E17').Value = '  -1.68%  '

# Row 18
$ws.Range('D18').Value = '67.044.66'
$ws.Range('E18').Value = '  -0.41%  '

# Row 19
$ws.Range('D19').Value = '''17.86'
$ws.Range('E19').Value = '  -3.13%  '

# Row 20
$ws.Range('D20').Value = '''11.65'
$ws.Range('E20').Value = '  -4.19%  '

# Row 21
$ws.Range('D21').Value = '''1.01'
$ws.Range('E21').Value = '  -6.02%  '

# Row 22
$ws.Range('D22').Value = '''393.10'
$ws.Range('E22').Value = '  -0.14%  '

# Row 23
$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').Value = '''3.91'
$ws.Range('E23').Value = '  -5.03%  '

# Row 24
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '''84.89'
$ws.Range('E24').Value = '  +0.89%  '

# Row 25
$ws.Range('D25').Value = '''11.52'
$ws.Range('E25').Value = '  -11.42%  '

# Row 26
$ws.Range('D26').Value = '''12.11'
$ws.Range('E26').Value = '  -1.67%  '

# Row 27
$ws.Range('D27').Value = '''2.78'
$ws.Range('E27').Value = '  -4.65%  '

# Row 28
$ws.Range('D28').Value = '''3.77'
$ws.Range('E28').Value = '  -0.61%  '

# Row 29
$ws.Range('D29').Value = '''8.69'
$ws.Range('E29').Value = '  -2.90%  '

# Row 30
$ws.Range('D30').Value = '''709.27'
$ws.Range('E30').Value = '  +4.57%  '

# Row 31
$ws.Range('D31').Value = '''30.64'
$ws.Range('E31').Value = '  -2.26%  '

# Row 32
$ws.Range('D32').Value = '''6.78'
$ws.Range('E32').Value = '  -17.45%  '

# Row 33
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').Value = '''63.69'
$ws.Range('E33').Value = '  +0.47%  '

# Row 34
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').Value = '''11.48'
$ws.Range('E34').Value = '  -4.91%  '

# Row 35
$ws.Range('D35').Value = '''0.110'
$ws.Range('E35').Value = '  -3.54%  '

# Row 36
$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D36').Value = '''1.00'
$ws.Range('E36').Value = '  -0.08%  '

# Row 37
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D37').Value = '''37.64'
$ws.Range('E37').Value = '  -10.84%  '

# Row 38
$ws.Range('D38').Value = '''0.383'
$ws.Range('E38').Value = '  -11.39%  '

# Row 39
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').Value = '''0.999'
$ws.Range('E39').Value = '  +0.27%  '

# Row 40
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '''0.129'
$ws.Range('E40').Value = '  -5.27%  '

# Row 41
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '3.029.31'
$ws.Range('E41').Value = '  -5.63%  '

# Row 42
$ws.Range('B42').Value = 'ThetaToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D42').Value = '''2.93'
$ws.Range('E42').Value = '  -6.15%  '

# Row 43
$ws.Range('D43').Value = '0.0₃0662'
$ws.Range('E43').Value = '  -13.41%  '

# Row 44
$ws.Range('E44').Value = '  +3.20%  '

# Row 45
$ws.Range('D45').Value = '''2.44'
$ws.Range('E45').Value = '  -14.04%  '

# Row 46
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '''0.0399'
$ws.Range('E46').Value = '  -3.22%  '

# Row 47
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '''0.130'
$ws.Range('E47').Value = '  +0.13%  '

# Row 48
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = '''2.93'
$ws.Range('E48').Value = '  -6.27%  '

# Row 49
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '''137.98'
$ws.Range('E49').Value = '  -1.21%  '

# Row 50
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '''2.62'
$ws.Range('E50').Value = '  -3.83%  '

# Row 51
$ws.Range('D51').Value = '''8.06'
$ws.Range('E51').Value = '  -5.95%  '
